{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"2025-07-21 Monday\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"2025-07-21 Monday\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"2025-07-22 Tuesday\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"597\u00d78=4776\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"597\u00d78=4776\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"327\u00d72=654\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"340\u00d77=2380\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"340\u00d77=2380\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"396\u00d73=1188\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"902\u00d75=4510\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"902\u00d75=4510\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"796\u00d74=3184\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"227\u00d74=908\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"227\u00d74=908\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"199\u00d74=796\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"463\u00d72=926\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"463\u00d72=926\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"709\u00d75=3545\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"635\u00d73=1905\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"635\u00d73=1905\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"241\u00d76=1446\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"126\u00d72=252\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"126\u00d72=252\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"839\u00d73=2517\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"267\u00d76=1602\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"267\u00d76=1602\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"992\u00d79=8928\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"256\u00d77=1792\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"256\u00d77=1792\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"296\u00d72=592\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"385\u00d73=1155\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"385\u00d73=1155\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"605\u00d73=1815\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"436\u00d79=3924\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"436\u00d79=3924\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"718\u00d78=5744\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"702\u00d79=6318\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"702\u00d79=6318\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"102\u00d72=204\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"125\u00d76=750\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"125\u00d76=750\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"898\u00d73=2694\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"463\u00d75=2315\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"463\u00d75=2315\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"357\u00d78=2856\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"246\u00d77=1722\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"246\u00d77=1722\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"658\u00d76=3948\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"662\u00d76=3972\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"662\u00d76=3972\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"160\u00d79=1440\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"302\u00d76=1812\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"302\u00d76=1812\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"909\u00d74=3636\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"648\u00d73=1944\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"648\u00d73=1944\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"855\u00d72=1710\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"186\u00d75=930\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"186\u00d75=930\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"116\u00d78=928\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"985\u00d78=7880\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"985\u00d78=7880\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"728\u00d73=2184\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"225\u00d72=450\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"225\u00d72=450\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"181\u00d74=724\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"405\u00d78=3240\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"405\u00d78=3240\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"816\u00d77=5712\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"767\u00d79=6903\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"767\u00d79=6903\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"747\u00d79=6723\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"902\u00d73=2706\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"902\u00d73=2706\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"105\u00d75=525\", Word.InsertLocation.replace);\n}\n\n{\n  const results = body.search(\"655\u00d73=1965\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \" + \"655\u00d73=1965\" + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"381\u00d72=762\", Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$found = $find.Execute(\"2025-07-21 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-07-22 Tuesday\", 2)\nif (-not $found) { throw \"Could not find text: 2025-07-21 Monday \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"597\u00d78=4776\", $false, $false, $false, $false, $false, $true, 1, $false, \"327\u00d72=654\", 2)\nif (-not $found) { throw \"Could not find text: 597\u00d78=4776 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"340\u00d77=2380\", $false, $false, $false, $false, $false, $true, 1, $false, \"396\u00d73=1188\", 2)\nif (-not $found) { throw \"Could not find text: 340\u00d77=2380 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"902\u00d75=4510\", $false, $false, $false, $false, $false, $true, 1, $false, \"796\u00d74=3184\", 2)\nif (-not $found) { throw \"Could not find text: 902\u00d75=4510 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"227\u00d74=908\", $false, $false, $false, $false, $false, $true, 1, $false, \"199\u00d74=796\", 2)\nif (-not $found) { throw \"Could not find text: 227\u00d74=908 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"463\u00d72=926\", $false, $false, $false, $false, $false, $true, 1, $false, \"709\u00d75=3545\", 2)\nif (-not $found) { throw \"Could not find text: 463\u00d72=926 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"635\u00d73=1905\", $false, $false, $false, $false, $false, $true, 1, $false, \"241\u00d76=1446\", 2)\nif (-not $found) { throw \"Could not find text: 635\u00d73=1905 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"126\u00d72=252\", $false, $false, $false, $false, $false, $true, 1, $false, \"839\u00d73=2517\", 2)\nif (-not $found) { throw \"Could not find text: 126\u00d72=252 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"267\u00d76=1602\", $false, $false, $false, $false, $false, $true, 1, $false, \"992\u00d79=8928\", 2)\nif (-not $found) { throw \"Could not find text: 267\u00d76=1602 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"256\u00d77=1792\", $false, $false, $false, $false, $false, $true, 1, $false, \"296\u00d72=592\", 2)\nif (-not $found) { throw \"Could not find text: 256\u00d77=1792 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"385\u00d73=1155\", $false, $false, $false, $false, $false, $true, 1, $false, \"605\u00d73=1815\", 2)\nif (-not $found) { throw \"Could not find text: 385\u00d73=1155 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"436\u00d79=3924\", $false, $false, $false, $false, $false, $true, 1, $false, \"718\u00d78=5744\", 2)\nif (-not $found) { throw \"Could not find text: 436\u00d79=3924 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"702\u00d79=6318\", $false, $false, $false, $false, $false, $true, 1, $false, \"102\u00d72=204\", 2)\nif (-not $found) { throw \"Could not find text: 702\u00d79=6318 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"125\u00d76=750\", $false, $false, $false, $false, $false, $true, 1, $false, \"898\u00d73=2694\", 2)\nif (-not $found) { throw \"Could not find text: 125\u00d76=750 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"463\u00d75=2315\", $false, $false, $false, $false, $false, $true, 1, $false, \"357\u00d78=2856\", 2)\nif (-not $found) { throw \"Could not find text: 463\u00d75=2315 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"246\u00d77=1722\", $false, $false, $false, $false, $false, $true, 1, $false, \"658\u00d76=3948\", 2)\nif (-not $found) { throw \"Could not find text: 246\u00d77=1722 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"662\u00d76=3972\", $false, $false, $false, $false, $false, $true, 1, $false, \"160\u00d79=1440\", 2)\nif (-not $found) { throw \"Could not find text: 662\u00d76=3972 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"302\u00d76=1812\", $false, $false, $false, $false, $false, $true, 1, $false, \"909\u00d74=3636\", 2)\nif (-not $found) { throw \"Could not find text: 302\u00d76=1812 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"648\u00d73=1944\", $false, $false, $false, $false, $false, $true, 1, $false, \"855\u00d72=1710\", 2)\nif (-not $found) { throw \"Could not find text: 648\u00d73=1944 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"186\u00d75=930\", $false, $false, $false, $false, $false, $true, 1, $false, \"116\u00d78=928\", 2)\nif (-not $found) { throw \"Could not find text: 186\u00d75=930 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"985\u00d78=7880\", $false, $false, $false, $false, $false, $true, 1, $false, \"728\u00d73=2184\", 2)\nif (-not $found) { throw \"Could not find text: 985\u00d78=7880 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"225\u00d72=450\", $false, $false, $false, $false, $false, $true, 1, $false, \"181\u00d74=724\", 2)\nif (-not $found) { throw \"Could not find text: 225\u00d72=450 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"405\u00d78=3240\", $false, $false, $false, $false, $false, $true, 1, $false, \"816\u00d77=5712\", 2)\nif (-not $found) { throw \"Could not find text: 405\u00d78=3240 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"767\u00d79=6903\", $false, $false, $false, $false, $false, $true, 1, $false, \"747\u00d79=6723\", 2)\nif (-not $found) { throw \"Could not find text: 767\u00d79=6903 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"902\u00d73=2706\", $false, $false, $false, $false, $false, $true, 1, $false, \"105\u00d75=525\", 2)\nif (-not $found) { throw \"Could not find text: 902\u00d73=2706 \" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"655\u00d73=1965\", $false, $false, $false, $false, $false, $true, 1, $false, \"381\u00d72=762\", 2)\nif (-not $found) { throw \"Could not find text: 655\u00d73=1965 \" }\n"}
